# Updates the cryptos price/volume table with freshly scraped values.
# Leading "'" on some D-column prices forces Excel to keep them as text
# (matching the source sheet, which stores every price/volume as a string)
# instead of silently parsing them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.772.40"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "2.045.30"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'227.70"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -2.14%  "
$ws.Range("D10").Value = "'0.0836"
$ws.Range("E10").Value = "  +2.95%  "
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").Value = "2.349.35"
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").Value = "'14.40"
$ws.Range("E13").Value = "  -1.65%  "
$ws.Range("D14").Value = "'21.42"
$ws.Range("E14").Value = "  +1.38%  "
$ws.Range("E15").Value = "  +6.27%  "
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "2.043.96"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "37.771.47"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("E20").Value = "  -1.88%  "
$ws.Range("D21").Value = "0.0₃0829"
$ws.Range("E21").Value = "  +0.74%  "
$ws.Range("D22").Value = "'222.49"
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("D24").Value = "'2.38"
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("E25").Value = "  +2.81%  "
$ws.Range("D26").Value = "'168.83"
$ws.Range("E26").Value = "  +2.26%  "
$ws.Range("D27").Value = "'9.32"
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("E28").Value = "  -1.28%  "
$ws.Range("E29").Value = "  -1.00%  "
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("E32").Value = "  +8.05%  "
$ws.Range("E33").Value = "  -1.46%  "
$ws.Range("D35").Value = "'0.0602"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "'6.48"
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("D37").Value = "'2.34"
$ws.Range("E37").Value = "  +3.90%  "
$ws.Range("D38").Value = "'3.47"
$ws.Range("E38").Value = "  +6.88%  "
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").Value = "'18.38"
$ws.Range("E40").Value = "  +9.46%  "
$ws.Range("D41").Value = "1.523.50"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").Value = "'97.54"
$ws.Range("E43").Value = "  -1.79%  "
$ws.Range("E44").Value = "  -0.14%  "
$ws.Range("D45").Value = "'4.18"
$ws.Range("E45").Value = "  +3.64%  "
$ws.Range("D46").Value = "'0.0891"
$ws.Range("E46").Value = "  -2.87%  "
$ws.Range("E47").Value = "  -0.35%  "
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "'7.10"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").Value = "'2.94"
$ws.Range("E50").Value = "  -0.39%  "
$ws.Range("D51").Value = "2.237.63"
$ws.Range("E51").Value = "  +0.77%  "
